$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.441.37'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '3.050.82'
$ws.Range("E3").Value = '  +4.40%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '201.91'
$ws.Range("D6").Value = '624.75'
$ws.Range("E6").Value = '  +4.65%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("D9").Value = '0.208'
$ws.Range("E9").Value = '  +5.15%  '
$ws.Range("D10").Value = '3.047.49'
$ws.Range("E10").Value = '  +4.36%  '
$ws.Range("D11").Value = "'0.440"
$ws.Range("E11").Value = '  +1.33%  '
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("D13").Value = '5.19'
$ws.Range("E13").Value = '  +5.80%  '
$ws.Range("D14").Value = '3.611.95'
$ws.Range("E14").Value = '  +4.51%  '
$ws.Range("D15").Value = '29.44'
$ws.Range("E15").Value = '  +4.52%  '
$ws.Range("D16").Value = '76.356.26'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("E17").Value = '  +2.16%  '
$ws.Range("D18").Value = '3.051.34'
$ws.Range("E18").Value = '  +4.64%  '
$ws.Range("D19").Value = '13.59'
$ws.Range("E19").Value = '  +4.86%  '
$ws.Range("D20").Value = '9.09'
$ws.Range("E20").Value = '  +4.15%  '
$ws.Range("D21").Value = '375.85'
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("D23").Value = '4.38'
$ws.Range("E23").Value = '  +2.34%  '
$ws.Range("D24").Value = '73.61'
$ws.Range("E24").Value = '  +3.11%  '
$ws.Range("D25").Value = '3.206.71'
$ws.Range("E25").Value = '  +4.50%  '
$ws.Range("E26").Value = '  +4.24%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  +2.30%  '
$ws.Range("E29").Value = '  +3.48%  '
$ws.Range("D30").Value = '0.996'
$ws.Range("E30").Value = '  -0.34%  '
$ws.Range("D31").Value = '8.33'
$ws.Range("E31").Value = '  +7.55%  '
$ws.Range("E32").Value = '  +1.17%  '
$ws.Range("D33").Value = '508.14'
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("E34").Value = '  +6.90%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '20.88'
$ws.Range("E36").Value = '  +3.20%  '
$ws.Range("D37").Value = '163.12'
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("E38").Value = '  +6.13%  '
$ws.Range("E39").Value = '  +2.11%  '
$ws.Range("D40").Value = '192.21'
$ws.Range("E40").Value = '  +4.85%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.116'
$ws.Range("E41").Value = '  +2.81%  '
$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").Value = '0.106'
$ws.Range("E42").Value = '  -5.39%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").Value = '5.19'
$ws.Range("E44").Value = '  +3.81%  '
$ws.Range("D45").Value = '0.787'
$ws.Range("E45").Value = '  +19.68%  '
$ws.Range("E46").Value = '  +8.19%  '
$ws.Range("D47").Value = '42.16'
$ws.Range("E47").Value = '  +5.15%  '
$ws.Range("E48").Value = '  +0.59%  '
$ws.Range("D49").Value = '2.49'
$ws.Range("E49").Value = '  +4.89%  '
$ws.Range("D50").Value = '0.614'
$ws.Range("E50").Value = '  +7.11%  '
$ws.Range("D51").Value = '3.96'
$ws.Range("E51").Value = '  +6.70%  '
